# UndoRedoActivityDiagram.pptx edit script
# - Update the cached "datetimeFigureOut" date placeholder text (6/7/2018 -> 14/4/19)
#   on every slide layout and on the slide master.
# - Rename "address book" -> "card collection" (and the matching camelCase
#   identifier) in the activity-diagram notes on slide 1.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Date placeholder ("6/7/2018" -> "14/4/19") on every layout + the master.
# ---------------------------------------------------------------------------
function Update-DateText($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "6/7/2018") {
                $sh.TextFrame.TextRange.Text = "14/4/19"
            }
        }
    }
}

$master = $p.SlideMaster
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateText($layouts.Item($li).Shapes)
}
Update-DateText($master.Shapes)

# ---------------------------------------------------------------------------
# 2. Slide 1 text: "address book" -> "card collection".
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if (-not $sh.HasTextFrame) {
        continue
    }
    $full = $sh.TextFrame.TextRange.Text

    if ($full -eq "[command commits address book]") {
        # Two runs: "[" and "command commits address book]" - only the
        # second run's text changes.
        $sh.TextFrame.TextRange.Runs(2).Text = "command commits card collection]"
    }

    if ($full -eq "Purge redundant states and then save address book to addressBookStateList ") {
        $tr = $sh.TextFrame.TextRange
        # Run 1 (chars 1-53): "Purge redundant states and then save address book to "
        $run1 = $tr.Characters(1, 53)
        $run1.Text = "Purge redundant states and then save card collection to "

        # Run 2 (chars 54-73 of the *new* text): "addressBookStateList"
        $tr2 = $sh.TextFrame.TextRange
        $prefixLen = "Purge redundant states and then save card collection to ".Length
        $run2 = $tr2.Characters($prefixLen + 1, 20)
        $run2.Text = "cardCollectionStateList"
        # Run 3 (trailing space) is left untouched.
    }
}
